$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Закупки"
$ws.Range("J4").Select()
